# Generate Report for Handoff
# Refresh the localization-status report: a new handoff run replaced the
# old source-doc GUID (70547e75-037f-404d-8dd5-9b72e7180881) with a new one
# (52451348-a402-4c40-b27b-006337d55dc7), produced new per-locale .xlf
# packages (new content hash), and re-stamped the handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldId = "70547e75-037f-404d-8dd5-9b72e7180881"
$newId = "52451348-a402-4c40-b27b-006337d55dc7"

$oldMd = "$oldId.md"
$newMd = "$newId.md"

$oldZhXlf = "$oldId.2f77d0cd2909cc630a82374c86b493b8ac9673e6.zh-cn.xlf"
$newZhXlf = "$newId.6f5293f97645b94a5b5efefa278da46f8456d2de.zh-cn.xlf"

$oldDeXlf = "$oldId.2f77d0cd2909cc630a82374c86b493b8ac9673e6.de-de.xlf"
$newDeXlf = "$newId.6f5293f97645b94a5b5efefa278da46f8456d2de.de-de.xlf"

$oldZhTime = "2016-03-10 13:01:31"
$newZhTime = "2016-03-10 13:01:51"

$oldDeTime = "2016-03-10 13:01:34"
$newDeTime = "2016-03-10 13:01:54"

function Update-DisplayedHyperlink($ws, $cellRef, $newText, $address) {
    # Update the text shown for an existing hyperlinked cell, keeping the
    # same external target. Setting TextToDisplay rewrites the cell's text
    # (and the hyperlink's display caption) without touching cell styles;
    # re-asserting Address afterwards keeps the link itself clickable/valid.
    # NOTE: index-based iteration (not foreach) is required here so the
    # Hyperlink handle stays "live" for the mutation below.
    $range = $ws.Range($cellRef)
    $hl = $null
    for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
        $candidate = $ws.Hyperlinks.Item($i)
        if ($candidate.Range.Address() -eq $range.Address()) {
            $hl = $candidate
        }
    }
    $hl.TextToDisplay = $newText
    $hl2 = $ws.Hyperlinks.Item($ws.Hyperlinks.Count)
    $hl2.Address = $address
}

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/1e139abf38be0b8b940e9ffd27972770bb90cee9/e2e/$oldMd"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff826a2fe09e76e74504f4fe08b401bcde76a47a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhXlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0b5a5fe8736552ff253f9800fe7b51fb72f5aec6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeXlf"

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Update-DisplayedHyperlink $wsOverview "A2" $newMd $mdAddress

# --- zh-cn sheet -------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-DisplayedHyperlink $wsZh "A2" $newMd $mdAddress
Update-DisplayedHyperlink $wsZh "C2" $newZhXlf $zhXlfAddress
$wsZh.Range("D2").Value = $newZhTime

# --- de-de sheet ---------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
Update-DisplayedHyperlink $wsDe "A2" $newMd $mdAddress
Update-DisplayedHyperlink $wsDe "C2" $newDeXlf $deXlfAddress
$wsDe.Range("D2").Value = $newDeTime
